$wb = $excel.ActiveWorkbook
$xlPasteValues = -4163

# ---------------------------------------------------------------------------
# Sheet 1: 台指期換倉成本計算  (A1:F6 -> A1:F7)
# Insert a new row 2 (pushing existing rows down) with the latest date entry.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("台指期換倉成本計算")
$ws1.Rows.Item(2).Insert()
$ws1.Range("A2").Value = "日期：2021/12/21"
# B2 ("202202") looks purely numeric - round-trip it through a formula +
# paste-values so it lands as literal text instead of being coerced to a
# number.
$ws1.Range("Z1").Formula = "=""202202"""
$ws1.Range("Z1").Copy()
$ws1.Range("B2").PasteSpecial($xlPasteValues)
$ws1.Range("Z1").Clear()
$ws1.Range("C2").Value = 17766
$ws1.Range("D2").Value = 2462
$ws1.Range("E2").Value = 6271398
$ws1.Range("F2").Value = 17636

# ---------------------------------------------------------------------------
# Sheet 2: 散戶多空力道  (A1:B40 -> A1:B41)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("散戶多空力道")
$ws2.Rows.Item(2).Insert()
$ws2.Range("A2").Value = "日期：2021/12/21"
$ws2.Range("B2").Value = 0.08

# ---------------------------------------------------------------------------
# Sheet 3: 三大法人買賣金額  (A1:C40 -> A1:C41)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("三大法人買賣金額")
$ws3.Rows.Item(2).Insert()
$ws3.Range("A2").Value = "110年12月21日"
$ws3.Range("B2").Value = 106.52
$ws3.Range("C2").Value = 13.64

# ---------------------------------------------------------------------------
# Sheet 4: 大盤多空點位  (A1:B39 -> A1:B40)
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("大盤多空點位")
$ws4.Rows.Item(2).Insert()
$ws4.Range("A2").Value = "110年12月21日"
$ws4.Range("B2").Value = 17748.15

# ---------------------------------------------------------------------------
# Sheet 5: 期貨大額交易人未沖銷部位  (A1:N38 -> A1:N39)
# ---------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("期貨大額交易人未沖銷部位")
$ws5.Rows.Item(2).Insert()
# A2 ("2021/12/21") looks like a date - round-trip it through a formula +
# paste-values so it lands as literal text instead of being coerced to a
# date serial.
$ws5.Range("Z1").Formula = "=""2021/12/21"""
$ws5.Range("Z1").Copy()
$ws5.Range("A2").PasteSpecial($xlPasteValues)
$ws5.Range("Z1").Clear()
$ws5.Range("B2").Value = 47439
$ws5.Range("C2").Value = 55318
$ws5.Range("D2").Value = 495
$ws5.Range("E2").Value = -78
$ws5.Range("F2").Value = 22106
$ws5.Range("G2").Value = 47519
$ws5.Range("H2").Value = 517
$ws5.Range("I2").Value = -712
$ws5.Range("J2").Value = -25413
$ws5.Range("K2").Value = 1229
$ws5.Range("L2").Value = -22
$ws5.Range("M2").Value = 634
$ws5.Range("N2").Value = -656
